$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Block: 150apps (rows 33-39), mirrors the existing 20/30/50/100 apps tables
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = "150apps"
$ws.Range("B33").Value = "HOM"
$ws.Range("C33").Value = "COM"
$ws.Range("D33").Value = "V-M"
$ws.Range("E33").Value = "ARI"
$ws.Range("F33").Value = "AMI"

$ws.Range("A34").Value = "Sift"
$ws.Range("B34").Value = 0.341
$ws.Range("C34").Value = 0.421
$ws.Range("D34").Value = 0.377
$ws.Range("E34").Value = 0.026
$ws.Range("F34").Value = 0.084

$ws.Range("A35").Value = "Hog"
$ws.Range("B35").Value = 0.355
$ws.Range("C35").Value = 0.424
$ws.Range("D35").Value = 0.387
$ws.Range("E35").Value = 0.026
$ws.Range("F35").Value = 0.076

$ws.Range("A36").Value = "Gabor"
$ws.Range("B36").Value = 0.338
$ws.Range("C36").Value = 0.419
$ws.Range("D36").Value = 0.374
$ws.Range("E36").Value = 0.028
$ws.Range("F36").Value = 0.08

$ws.Range("A37").Value = "AE"
$ws.Range("B37").Value = 0.455
$ws.Range("C37").Value = 0.53
$ws.Range("D37").Value = 0.49
$ws.Range("E37").Value = 0.084
$ws.Range("F37").Value = 0.198

$ws.Range("A38").Value = "GUI2Vec"
$ws.Range("B38").Value = 0.464
$ws.Range("C38").Value = 0.558
$ws.Range("D38").Value = 0.507
$ws.Range("E38").Value = 0.104
$ws.Range("F38").Value = 0.236
$ws.Range("B38:F38").Font.Color = 12611584

$ws.Range("A39").Value = "improvement"
$ws.Range("B39").Formula = "=(B38-B37)/B37"
$ws.Range("C39").Formula = "=(C38-C37)/C37"
$ws.Range("D39").Formula = "=(D38-D37)/D37"
$ws.Range("E39").Formula = "=(E38-E37)/E37"
$ws.Range("F39").Formula = "=(F38-F37)/F37"

# ---------------------------------------------------------------------------
# Block: 200apps (rows 41-47)
# ---------------------------------------------------------------------------
$ws.Range("A41").Value = "200apps"
$ws.Range("B41").Value = "HOM"
$ws.Range("C41").Value = "COM"
$ws.Range("D41").Value = "V-M"
$ws.Range("E41").Value = "ARI"
$ws.Range("F41").Value = "AMI"

$ws.Range("A42").Value = "Sift"
$ws.Range("B42").Value = 0.347
$ws.Range("C42").Value = 0.438
$ws.Range("D42").Value = 0.387
$ws.Range("E42").Value = 0.022
$ws.Range("F42").Value = 0.079

$ws.Range("A43").Value = "Hog"
$ws.Range("B43").Value = 0.374
$ws.Range("C43").Value = 0.444
$ws.Range("D43").Value = 0.406
$ws.Range("E43").Value = 0.02
$ws.Range("F43").Value = 0.069

$ws.Range("A44").Value = "Gabor"
$ws.Range("B44").Value = 0.369
$ws.Range("C44").Value = 0.444
$ws.Range("D44").Value = 0.403
$ws.Range("E44").Value = 0.026
$ws.Range("F44").Value = 0.074

$ws.Range("A45").Value = "AE"
$ws.Range("B45").Value = 0.45
$ws.Range("C45").Value = 0.525
$ws.Range("D45").Value = 0.485
$ws.Range("E45").Value = 0.064
$ws.Range("F45").Value = 0.163

$ws.Range("A46").Value = "GUI2Vec"
$ws.Range("B46").Value = 0.453
$ws.Range("C46").Value = 0.555
$ws.Range("D46").Value = 0.499
$ws.Range("E46").Value = 0.076
$ws.Range("F46").Value = 0.207

$ws.Range("A47").Value = "improvement"
$ws.Range("B47").Formula = "=(B46-B45)/B45"
$ws.Range("C47").Formula = "=(C46-C45)/C45"
$ws.Range("D47").Formula = "=(D46-D45)/D45"
$ws.Range("E47").Formula = "=(E46-E45)/E45"
$ws.Range("F47").Formula = "=(F46-F45)/F45"

# ---------------------------------------------------------------------------
# Sheet view: drop the old topLeftCell/selection, select F3 instead
# ---------------------------------------------------------------------------
$null = $ws.Range("F3").Select()
